$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 1) Title: "Bulls and Cows console game" -> 3 runs, "Cows" -> "Goats"
$xml = '<w:p w:rsidR="004539CB" w:rsidRDefault="00A30854" w:rsidP="00A4202E"><w:pPr><w:pStyle w:val="Title"/></w:pPr>'
$xml += '<w:r w:rsidRPr="00A30854"><w:t xml:space="preserve">Bulls and </w:t></w:r>'
$xml += '<w:r w:rsidRPr="00A30854"><w:t>Goats</w:t></w:r>'
$xml += '<w:r w:rsidRPr="00A30854"><w:t xml:space="preserve"> console game</w:t></w:r>'
$xml += '</w:p>'
$null = $d.Paragraphs(2).Range.InsertXML($pkgHeader + $xml + $pkgFooter)

# 2) "The gamer" + " want" + "s:" -> single run "The gamer wants:"
$xml = '<w:p w:rsidR="00A30854" w:rsidRDefault="00A30854" w:rsidP="005D3A98"><w:pPr><w:spacing w:after="0"/></w:pPr>'
$xml += '<w:r><w:t>The gamer wants:</w:t></w:r>'
$xml += '</w:p>'
$null = $d.Paragraphs(5).Range.InsertXML($pkgHeader + $xml + $pkgFooter)

# 3) "A guessing game with words" -> append new run " based on Mastermind"
$xml = '<w:p w:rsidR="00A30854" w:rsidRDefault="00A30854" w:rsidP="00A30854"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>'
$xml += '<w:r><w:t>A guessing game with words</w:t></w:r>'
$xml += '<w:r><w:t xml:space="preserve"> based on Mastermind</w:t></w:r>'
$xml += '</w:p>'
$null = $d.Paragraphs(14).Range.InsertXML($pkgHeader + $xml + $pkgFooter)

# 4) "Isogram" + "s" + ": the words have no repeating letters" -> single run
$xml = '<w:p w:rsidR="00A30854" w:rsidRDefault="00A30854" w:rsidP="00A30854"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>'
$xml += '<w:r><w:t>Isograms: the words have no repeating letters</w:t></w:r>'
$xml += '</w:p>'
$null = $d.Paragraphs(15).Range.InsertXML($pkgHeader + $xml + $pkgFooter)

# 5) "Cow:  the no. ..." + " (... umm cow's eye?? " -> "Goat" / ":  the no. ..." /
#    " (... umm " / "goat" / "'s eye?? ", keeping the existing <w:sym/> + " )" runs
$ellipsis = [char]0x2026
$rsquo = [char]0x2019
$xml = '<w:p w:rsidR="00A30854" w:rsidRDefault="00A30854" w:rsidP="00A30854"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>'
$xml += '<w:r><w:t>Goat</w:t></w:r>'
$xml += '<w:r><w:t>:  the no. of right letters in the wrong place</w:t></w:r>'
$xml += '<w:r w:rsidR="005D3A98"><w:t xml:space="preserve"> (' + $ellipsis + ' from ' + $ellipsis + ' umm </w:t></w:r>'
$xml += '<w:r w:rsidR="005D3A98"><w:t>goat</w:t></w:r>'
$xml += '<w:r w:rsidR="005D3A98"><w:t>' + $rsquo + 's eye?? </w:t></w:r>'
$xml += '<w:r w:rsidR="005D3A98"><w:sym w:font="Wingdings" w:char="F04A"/></w:r>'
$xml += '<w:r w:rsidR="005D3A98"><w:t xml:space="preserve"> )</w:t></w:r>'
$xml += '</w:p>'
$null = $d.Paragraphs(20).Range.InsertXML($pkgHeader + $xml + $pkgFooter)

# 6) "no" + ". " + "of bulls and cows" -> "no. " / "of bulls and " / "goats"
$xml = '<w:p w:rsidR="00A30854" w:rsidRDefault="00A30854" w:rsidP="008C1ADA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>'
$xml += '<w:r><w:t xml:space="preserve">the </w:t></w:r>'
$xml += '<w:r w:rsidR="00122ED4"><w:t xml:space="preserve">no. </w:t></w:r>'
$xml += '<w:r><w:t xml:space="preserve">of bulls and </w:t></w:r>'
$xml += '<w:r><w:t>goats</w:t></w:r>'
$xml += '</w:p>'
$null = $d.Paragraphs(30).Range.InsertXML($pkgHeader + $xml + $pkgFooter)

# 7) "ASCII art of a bull and a cow" -> split into 2 runs, "cow" -> "goat"
$xml = '<w:p w:rsidR="005D3A98" w:rsidRDefault="005D3A98" w:rsidP="005D3A98"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>'
$xml += '<w:r><w:t xml:space="preserve">ASCII art of a bull and a </w:t></w:r>'
$xml += '<w:r><w:t>goat</w:t></w:r>'
$xml += '</w:p>'
$null = $d.Paragraphs(45).Range.InsertXML($pkgHeader + $xml + $pkgFooter)

Write-Output "Done"
